# Add MON (9/6/2021, Labor Day) as a shaded/holiday column on the "1-15"
# sign-in sheet, matching the styling already used for the SAT/SUN columns
# (columns H/I), and narrow columns L/M from 4.5 to the same 2.5 width as
# the other holiday columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1-15")

# --- Column widths: L (12) and M (13) go from 4.5 -> 2.5, matching H/I ---
$holidayWidth = $ws.Columns.Item(8).ColumnWidth
$ws.Columns.Item(12).ColumnWidth = $holidayWidth
$ws.Columns.Item(13).ColumnWidth = $holidayWidth

# Rows whose L/M cells take on the SAT/SUN-style gray-fill/border formatting.
# Rows with an "X" entry (clinic rows) also get an "X" written into L & M.
$rowsWithX = @(5,6,8,9,11,12,14,15,17,18,20,21,23,24,26,27)
$allRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27)

foreach ($r in $allRows) {
    $ws.Range("H$r").Copy()
    $ws.Range("L$r").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("I$r").Copy()
    $ws.Range("M$r").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

foreach ($r in $rowsWithX) {
    $ws.Range("L$r").Value = "X"
    $ws.Range("M$r").Value = "X"
}
